$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.044816
$ws.Range("H2").Value = 0.134448
$ws.Range("I2").Value = 0.031000309200692
$ws.Range("J2").Value = 0.031000309200692
$ws.Range("Q2").Value = 0.5161562394346667
$ws.Range("R2").Value = 4.645406154912
$ws.Range("S2").Value = 0.031000309200692
$ws.Range("T2").Value = 0.031000309200692

# Row 3 updates
$ws.Range("I3").Value = 0.5689052935112355
$ws.Range("J3").Value = 0.5689052935112355
$ws.Range("S3").Value = 0.5689052935112355
$ws.Range("T3").Value = 0.5689052935112355

# Row 4 updates
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5784016666666667
$ws.Range("H4").Value = 1.735205
$ws.Range("I4").Value = 0.4000943972880724
$ws.Range("J4").Value = 0.4000943972880724
$ws.Range("Q4").Value = 6.661585798585556
$ws.Range("R4").Value = 59.95427218727
$ws.Range("S4").Value = 0.4000943972880724
$ws.Range("T4").Value = 0.4000943972880724
